$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each row based on the latest scrape.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.809.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.076.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.50%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.27%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.02"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.92%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0788"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.78%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.79"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.26%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.774"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.36"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.096.31"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "37.707.62"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.85%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.54"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0850"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.25"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.42"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.92%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.92"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.67%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.63%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.50"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.71"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.77"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0633"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.17%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.82"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0975"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.30"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.41%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.59"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.440.84"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.18%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.17"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.78%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.43%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.53%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.267.43"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.68"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.33%  "
